$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 2441.9
$ws.Range("I11").Value = 2441.9
$ws.Range("K11").Value = 2441.9
$ws.Range("M11").Value = -2301.9
$ws.Range("H40").Value = 3621.4614
$ws.Range("I40").Value = 3211.111
$ws.Range("J40").Value = 4544.75
$ws.Range("K40").Value = 3211.111
$ws.Range("L40").Value = 4544.75
$ws.Range("M40").Value = -3036.111
$ws.Range("N40").Value = -4894.75
$ws.Range("H76").Value = 7696381
$ws.Range("J76").Value = 4328
$ws.Range("L76").Value = 4328
$ws.Range("N76").Value = -4958
$ws.Range("H79").Value = 7696381
$ws.Range("J79").Value = 4328
$ws.Range("L79").Value = 4328
$ws.Range("N79").Value = -6512
$ws.Range("H98").Value = 47621550
$ws.Range("I98").Value = 55558308
$ws.Range("K98").Value = 55558308
$ws.Range("M98").Value = -55556810
$ws.Range("H100").Value = 2618.9443
$ws.Range("I100").Value = 1743
$ws.Range("J100").Value = 3995.4285
$ws.Range("K100").Value = 1743
$ws.Range("L100").Value = 3995.4285
$ws.Range("M100").Value = -1202
$ws.Range("N100").Value = -5077.4285
$ws.Range("H122").Value = 47621550
$ws.Range("I122").Value = 55558308
$ws.Range("K122").Value = 166674924
$ws.Range("M122").Value = -166672474
$ws.Range("H132").Value = 1122.5
$ws.Range("I132").Value = 1074.5918
$ws.Range("K132").Value = 3223.7754
$ws.Range("M132").Value = -693.7753999999995
$ws.Range("H135").Value = 294714.6
$ws.Range("I135").Value = 313109.25
$ws.Range("K135").Value = 2817983.25
$ws.Range("M135").Value = -2815448.25
$ws.Range("H138").Value = 3570.3535
$ws.Range("I138").Value = 962.4103
$ws.Range("J138").Value = 5265.5166
$ws.Range("K138").Value = 2887.2309
$ws.Range("L138").Value = 15796.5498
$ws.Range("M138").Value = 2252.7691
$ws.Range("N138").Value = -26076.5498
$ws.Range("H141").Value = 8334761.5
$ws.Range("J141").Value = 3095.3333
$ws.Range("L141").Value = 9285.999899999999
$ws.Range("N141").Value = -19645.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 933.8
$ws.Range("I32").Value = 934.6836499999999
$ws.Range("K32").Value = 934.6836499999999
$ws.Range("M32").Value = -647.6836499999999
$ws.Range("H45").Value = 8478.4
$ws.Range("I45").Value = 6969.143
$ws.Range("K45").Value = 6969.143
$ws.Range("M45").Value = -6592.143
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H74").Value = 26588.45
$ws.Range("I74").Value = 34237.535
$ws.Range("K74").Value = 34237.535
$ws.Range("M74").Value = -33363.535
$ws.Range("H77").Value = 26588.45
$ws.Range("I77").Value = 34237.535
$ws.Range("K77").Value = 171187.675
$ws.Range("M77").Value = -166819.675
$ws.Range("H97").Value = 11905032
$ws.Range("I97").Value = 247.8
$ws.Range("J97").Value = 41666990
$ws.Range("K97").Value = 247.8
$ws.Range("L97").Value = 41666990
$ws.Range("M97").Value = 248.2
$ws.Range("N97").Value = -41667982

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3832.5
$ws.Range("J94").Value = 7439.875
$ws.Range("L94").Value = 7439.875
$ws.Range("N94").Value = -8341.875
$ws.Range("H105").Value = 2633.8845
$ws.Range("I105").Value = 1752.2667
$ws.Range("J105").Value = 3836.0908
$ws.Range("K105").Value = 1752.2667
$ws.Range("L105").Value = 3836.0908
$ws.Range("M105").Value = -5.266699999999901
$ws.Range("N105").Value = -7330.0908
$ws.Range("H107").Value = 80361090
$ws.Range("I107").Value = 225000690
$ws.Range("K107").Value = 225000690
$ws.Range("M107").Value = -224998770
$ws.Range("H123").Value = 51598
$ws.Range("J123").Value = 51598
$ws.Range("L123").Value = 51598
$ws.Range("N123").Value = -61398
$ws.Range("H134").Value = 6072.841
$ws.Range("I134").Value = 2013.4
$ws.Range("K134").Value = 6040.200000000001
$ws.Range("M134").Value = -3505.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6011.9106
$ws.Range("I31").Value = 2557.4146
$ws.Range("K31").Value = 2557.4146
$ws.Range("M31").Value = -2262.4146
$ws.Range("H34").Value = 6011.9106
$ws.Range("I34").Value = 2557.4146
$ws.Range("K34").Value = 2557.4146
$ws.Range("M34").Value = -2355.4146
$ws.Range("H42").Value = 42833.332
$ws.Range("J42").Value = 42833.332
$ws.Range("L42").Value = 42833.332
$ws.Range("N42").Value = -44019.332
$ws.Range("H43").Value = 32656.5
$ws.Range("J43").Value = 32656.5
$ws.Range("L43").Value = 32656.5
$ws.Range("N43").Value = -33024.5
$ws.Range("H58").Value = 7466805.5
$ws.Range("I58").Value = 11365045
$ws.Range("K58").Value = 11365045
$ws.Range("M58").Value = -11364842
$ws.Range("H99").Value = 7734.6665
$ws.Range("I99").Value = 3799
$ws.Range("K99").Value = 3799
$ws.Range("M99").Value = -2301
$ws.Range("H101").Value = 32656.5
$ws.Range("J101").Value = 32656.5
$ws.Range("L101").Value = 32656.5
$ws.Range("N101").Value = -39146.5
$ws.Range("H102").Value = 44773.8
$ws.Range("J102").Value = 44773.8
$ws.Range("L102").Value = 44773.8
$ws.Range("N102").Value = -49641.8
$ws.Range("H107").Value = 2175.6843
$ws.Range("I107").Value = 572.1667
$ws.Range("J107").Value = 2915.7693
$ws.Range("K107").Value = 572.1667
$ws.Range("L107").Value = 2915.7693
$ws.Range("M107").Value = 1347.8333
$ws.Range("N107").Value = -6755.7693
$ws.Range("H126").Value = 7734.6665
$ws.Range("I126").Value = 3799
$ws.Range("K126").Value = 11397
$ws.Range("M126").Value = -8927
$ws.Range("H132").Value = 3611.5972
$ws.Range("I132").Value = 1557.8113
$ws.Range("K132").Value = 4673.4339
$ws.Range("M132").Value = -2143.4339
$ws.Range("H134").Value = 4247.0293
$ws.Range("I134").Value = 1771.4348
$ws.Range("K134").Value = 5314.3044
$ws.Range("M134").Value = -2779.3044
$ws.Range("H136").Value = 7466805.5
$ws.Range("I136").Value = 11365045
$ws.Range("K136").Value = 34095135
$ws.Range("M136").Value = -34092585

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1381603.2
$ws.Range("I5").Value = 2500612.5
$ws.Range("J5").Value = 4361.3076
$ws.Range("K5").Value = 7501837.5
$ws.Range("L5").Value = 13083.9228
$ws.Range("M5").Value = -7501725.5
$ws.Range("N5").Value = -13307.9228
$ws.Range("H44").Value = 1177.1428
$ws.Range("J44").Value = 1812.5
$ws.Range("L44").Value = 5437.5
$ws.Range("N44").Value = -6233.5
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H97").Value = 399.5
$ws.Range("J97").Value = 399.5
$ws.Range("L97").Value = 1198.5
$ws.Range("N97").Value = -2190.5
$ws.Range("H129").Value = 68203
$ws.Range("J129").Value = 78460.38
$ws.Range("L129").Value = 235381.14
$ws.Range("N129").Value = -245381.14
$ws.Range("H131").Value = 1648.45
$ws.Range("I131").Value = 1414.25
$ws.Range("J131").Value = 1804.5834
$ws.Range("K131").Value = 4242.75
$ws.Range("L131").Value = 5413.7502
$ws.Range("M131").Value = 797.25
$ws.Range("N131").Value = -15493.7502
$ws.Range("H135").Value = 1381603.2
$ws.Range("I135").Value = 2500612.5
$ws.Range("J135").Value = 4361.3076
$ws.Range("K135").Value = 22505512.5
$ws.Range("L135").Value = 39251.7684
$ws.Range("M135").Value = -22502977.5
$ws.Range("N135").Value = -44321.7684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1072.2858
$ws.Range("I13").Value = 302
$ws.Range("J13").Value = 2998
$ws.Range("K13").Value = 302
$ws.Range("L13").Value = 2998
$ws.Range("M13").Value = -163
$ws.Range("N13").Value = -3276
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 30245
$ws.Range("J32").Value = 30245
$ws.Range("L32").Value = 30245
$ws.Range("N32").Value = -30837
$ws.Range("H80").Value = 2768.2856
$ws.Range("J80").Value = 2300
$ws.Range("L80").Value = 2300
$ws.Range("N80").Value = -4296
$ws.Range("H83").Value = 2768.2856
$ws.Range("J83").Value = 2300
$ws.Range("L83").Value = 11500
$ws.Range("N83").Value = -21484
$ws.Range("H122").Value = 3151036.8
$ws.Range("I122").Value = 4261756
$ws.Range("J122").Value = 3999.3333
$ws.Range("K122").Value = 12785268
$ws.Range("L122").Value = 11997.9999
$ws.Range("M122").Value = -12782818
$ws.Range("N122").Value = -16897.9999
$ws.Range("H132").Value = 7508.7
$ws.Range("I132").Value = 2595.8333
$ws.Range("J132").Value = 14878
$ws.Range("K132").Value = 7787.499899999999
$ws.Range("L132").Value = 44634
$ws.Range("M132").Value = -5257.499899999999
$ws.Range("N132").Value = -49694
$ws.Range("H135").Value = 78556.5
$ws.Range("J135").Value = 78556.5
$ws.Range("L135").Value = 78556.5
$ws.Range("N135").Value = -88696.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 15876407
$ws.Range("I46").Value = 2501
$ws.Range("J46").Value = 18522058
$ws.Range("K46").Value = 2501
$ws.Range("L46").Value = 18522058
$ws.Range("M46").Value = -2313
$ws.Range("N46").Value = -18522434
$ws.Range("H132").Value = 12202474
$ws.Range("I132").Value = 27781194
$ws.Range("K132").Value = 83343582
$ws.Range("M132").Value = -83341052
$ws.Range("H136").Value = 5813.1274
$ws.Range("I136").Value = 835.2414
$ws.Range("K136").Value = 2505.7242
$ws.Range("M136").Value = 44.27579999999989

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 56546
$ws.Range("J16").Value = 56546
$ws.Range("L16").Value = 56546
$ws.Range("N16").Value = -57130
$ws.Range("H116").Value = 55864.5
$ws.Range("J116").Value = 55864.5
$ws.Range("L116").Value = 55864.5
$ws.Range("N116").Value = -65042.5
$ws.Range("H117").Value = 54183.5
$ws.Range("J117").Value = 54183.5
$ws.Range("L117").Value = 54183.5
$ws.Range("N117").Value = -63361.5
$ws.Range("H119").Value = 56958
$ws.Range("J119").Value = 56958
$ws.Range("L119").Value = 56958
$ws.Range("N119").Value = -66634
$ws.Range("H120").Value = 56958
$ws.Range("J120").Value = 56958
$ws.Range("L120").Value = 56958
$ws.Range("N120").Value = -66634
$ws.Range("H121").Value = 54189
$ws.Range("J121").Value = 54189
$ws.Range("L121").Value = 54189
$ws.Range("N121").Value = -57683
$ws.Range("H123").Value = 47274.25
$ws.Range("J123").Value = 47274.25
$ws.Range("L123").Value = 47274.25
$ws.Range("N123").Value = -57074.25
$ws.Range("H124").Value = 46532.332
$ws.Range("J124").Value = 46532.332
$ws.Range("L124").Value = 46532.332
$ws.Range("N124").Value = -56352.332
$ws.Range("H125").Value = 84536.25
$ws.Range("J125").Value = 84536.25
$ws.Range("L125").Value = 84536.25
$ws.Range("N125").Value = -94376.25
$ws.Range("H129").Value = 50000
$ws.Range("J129").Value = 50000
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000
$ws.Range("H132").Value = 12510814
$ws.Range("I132").Value = 17246296
$ws.Range("K132").Value = 51738888
$ws.Range("M132").Value = -51736358
